# Apply scheduled-runner market-data refresh to the Phoenix profits workbook.
# Values come from the authoritative diff; this only rewrites cells in H:N
# (currentAveragePrice*, Leve price/profit columns) -- no formulas anywhere
# in this workbook, so plain value assignment reproduces the diff exactly.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item(1)
# row 51
$ws.Cells.Item(51,8).Value = 4721.35
$ws.Cells.Item(51,10).Value = 4491.4546
$ws.Cells.Item(51,12).Value = 4491.4546
$ws.Cells.Item(51,14).Value = -5459.4546

# row 64
$ws.Cells.Item(64,8).Value = 11347
$ws.Cells.Item(64,9).Value = 3701
$ws.Cells.Item(64,10).Value = 16444.334
$ws.Cells.Item(64,11).Value = 3701
$ws.Cells.Item(64,12).Value = 16444.334
$ws.Cells.Item(64,13).Value = -3453
$ws.Cells.Item(64,14).Value = -16940.334

# row 67
$ws.Cells.Item(67,8).Value = 11347
$ws.Cells.Item(67,9).Value = 3701
$ws.Cells.Item(67,10).Value = 16444.334
$ws.Cells.Item(67,11).Value = 3701
$ws.Cells.Item(67,12).Value = 16444.334
$ws.Cells.Item(67,13).Value = -2843
$ws.Cells.Item(67,14).Value = -18160.334

# row 70
$ws.Cells.Item(70,8).Value = 3071.2307
$ws.Cells.Item(70,9).Value = 1000
$ws.Cells.Item(70,10).Value = 3243.8333
$ws.Cells.Item(70,11).Value = 3000
$ws.Cells.Item(70,12).Value = 9731.499899999999
$ws.Cells.Item(70,13).Value = -2730
$ws.Cells.Item(70,14).Value = -10271.4999

# row 73
$ws.Cells.Item(73,8).Value = 3071.2307
$ws.Cells.Item(73,9).Value = 1000
$ws.Cells.Item(73,10).Value = 3243.8333
$ws.Cells.Item(73,11).Value = 3000
$ws.Cells.Item(73,12).Value = 9731.499899999999
$ws.Cells.Item(73,13).Value = -2064
$ws.Cells.Item(73,14).Value = -11603.4999

# row 74
$ws.Cells.Item(74,8).Value = 5491.6665
$ws.Cells.Item(74,9).Value = 5491.6665
$ws.Cells.Item(74,11).Value = 5491.6665
$ws.Cells.Item(74,13).Value = -4555.6665

# row 77
$ws.Cells.Item(77,8).Value = 5491.6665
$ws.Cells.Item(77,9).Value = 5491.6665
$ws.Cells.Item(77,11).Value = 27458.3325
$ws.Cells.Item(77,13).Value = -22778.3325

# row 132
$ws.Cells.Item(132,8).Value = 4563.2354
$ws.Cells.Item(132,9).Value = 4563.2354
$ws.Cells.Item(132,11).Value = 13689.7062
$ws.Cells.Item(132,13).Value = -11159.7062

# row 133
$ws.Cells.Item(133,8).Value = 190591
$ws.Cells.Item(133,10).Value = 190591
$ws.Cells.Item(133,12).Value = 190591
$ws.Cells.Item(133,14).Value = -200711

# row 137
$ws.Cells.Item(137,8).Value = 123356.22
$ws.Cells.Item(137,9).Value = 1435.8462
$ws.Cells.Item(137,10).Value = 232664.14
$ws.Cells.Item(137,11).Value = 4307.5386
$ws.Cells.Item(137,12).Value = 697992.42
$ws.Cells.Item(137,13).Value = -1757.5386
$ws.Cells.Item(137,14).Value = -703092.42

# row 139
$ws.Cells.Item(139,8).Value = 99923.5
$ws.Cells.Item(139,10).Value = 99923.5
$ws.Cells.Item(139,12).Value = 99923.5
$ws.Cells.Item(139,14).Value = -110203.5

# row 140
$ws.Cells.Item(140,8).Value = 75624.164
$ws.Cells.Item(140,10).Value = 75624.164
$ws.Cells.Item(140,12).Value = 75624.164
$ws.Cells.Item(140,14).Value = -85984.164

# row 141
$ws.Cells.Item(141,8).Value = 5758.706
$ws.Cells.Item(141,9).Value = 5099.8936
$ws.Cells.Item(141,11).Value = 15299.6808
$ws.Cells.Item(141,13).Value = -10119.6808

# --- ARM ---
$ws = $wb.Worksheets.Item(2)
# row 32
$ws.Cells.Item(32,8).Value = 4459.1304
$ws.Cells.Item(32,9).Value = 3929.439
$ws.Cells.Item(32,11).Value = 3929.439
$ws.Cells.Item(32,13).Value = -3642.439

# row 76
$ws.Cells.Item(76,8).Value = 500000
$ws.Cells.Item(76,10).Value = 500000
$ws.Cells.Item(76,12).Value = 500000
$ws.Cells.Item(76,14).Value = -500676

# row 79
$ws.Cells.Item(79,8).Value = 500000
$ws.Cells.Item(79,10).Value = 500000
$ws.Cells.Item(79,12).Value = 500000
$ws.Cells.Item(79,14).Value = -502340

# row 109
$ws.Cells.Item(109,8).Value = 53482.668
$ws.Cells.Item(109,10).Value = 53482.668
$ws.Cells.Item(109,12).Value = 53482.668
$ws.Cells.Item(109,14).Value = -56256.668

# row 132
$ws.Cells.Item(132,8).Value = 5566.471
$ws.Cells.Item(132,9).Value = 5634.508
$ws.Cells.Item(132,11).Value = 16903.524
$ws.Cells.Item(132,13).Value = -14373.524

# --- BSM ---
$ws = $wb.Worksheets.Item(3)
# row 86
$ws.Cells.Item(86,8).Value = 28574042
$ws.Cells.Item(86,9).Value = 47621820
$ws.Cells.Item(86,10).Value = 2377.7856
$ws.Cells.Item(86,11).Value = 47621820
$ws.Cells.Item(86,12).Value = 2377.7856
$ws.Cells.Item(86,13).Value = -47620697
$ws.Cells.Item(86,14).Value = -4623.7856

# row 89
$ws.Cells.Item(89,8).Value = 28574042
$ws.Cells.Item(89,9).Value = 47621820
$ws.Cells.Item(89,10).Value = 2377.7856
$ws.Cells.Item(89,11).Value = 238109100
$ws.Cells.Item(89,12).Value = 11888.928
$ws.Cells.Item(89,13).Value = -238103484
$ws.Cells.Item(89,14).Value = -23120.928

# row 138
$ws.Cells.Item(138,8).Value = 98414.375
$ws.Cells.Item(138,10).Value = 98414.375
$ws.Cells.Item(138,12).Value = 98414.375
$ws.Cells.Item(138,14).Value = -108694.375

# --- CRP ---
$ws = $wb.Worksheets.Item(4)
# row 31
$ws.Cells.Item(31,8).Value = 2095.5386
$ws.Cells.Item(31,9).Value = 2001.7
$ws.Cells.Item(31,11).Value = 2001.7
$ws.Cells.Item(31,13).Value = -1706.7

# row 34
$ws.Cells.Item(34,8).Value = 2095.5386
$ws.Cells.Item(34,9).Value = 2001.7
$ws.Cells.Item(34,11).Value = 2001.7
$ws.Cells.Item(34,13).Value = -1799.7

# row 86
$ws.Cells.Item(86,8).Value = 17998.285
$ws.Cells.Item(86,10).Value = 19331.334
$ws.Cells.Item(86,12).Value = 19331.334
$ws.Cells.Item(86,14).Value = -21577.334

# row 89
$ws.Cells.Item(89,8).Value = 17998.285
$ws.Cells.Item(89,10).Value = 19331.334
$ws.Cells.Item(89,12).Value = 96656.67
$ws.Cells.Item(89,14).Value = -107888.67

# row 102
$ws.Cells.Item(102,8).Value = 37299
$ws.Cells.Item(102,10).Value = 38449.5
$ws.Cells.Item(102,12).Value = 38449.5
$ws.Cells.Item(102,14).Value = -43317.5

# row 105
$ws.Cells.Item(105,8).Value = 2534.8
$ws.Cells.Item(105,9).Value = 2534.8
$ws.Cells.Item(105,11).Value = 2534.8
$ws.Cells.Item(105,13).Value = -787.8000000000002

# row 115
$ws.Cells.Item(115,8).Value = 37633
$ws.Cells.Item(115,10).Value = 37633
$ws.Cells.Item(115,12).Value = 37633
$ws.Cells.Item(115,14).Value = -39983

# row 134
$ws.Cells.Item(134,8).Value = 4173.452
$ws.Cells.Item(134,9).Value = 4231.9375
$ws.Cells.Item(134,10).Value = 3986.3
$ws.Cells.Item(134,11).Value = 12695.8125
$ws.Cells.Item(134,12).Value = 11958.9
$ws.Cells.Item(134,13).Value = -10160.8125
$ws.Cells.Item(134,14).Value = -17028.9

# row 138
$ws.Cells.Item(138,8).Value = 86111.14
$ws.Cells.Item(138,10).Value = 86111.14
$ws.Cells.Item(138,12).Value = 86111.14
$ws.Cells.Item(138,14).Value = -96391.14

# --- CUL ---
$ws = $wb.Worksheets.Item(5)
# row 68
$ws.Cells.Item(68,8).Value = 1014.8
$ws.Cells.Item(68,9).Value = 946.5833
$ws.Cells.Item(68,11).Value = 2839.7499
$ws.Cells.Item(68,13).Value = -2028.7499

# row 71
$ws.Cells.Item(71,8).Value = 1014.8
$ws.Cells.Item(71,9).Value = 946.5833
$ws.Cells.Item(71,11).Value = 8519.2497
$ws.Cells.Item(71,13).Value = -4463.2497

# --- GSM ---
$ws = $wb.Worksheets.Item(6)
# row 3
$ws.Cells.Item(3,8).Value = 1884.7142
$ws.Cells.Item(3,9).Value = 3075
$ws.Cells.Item(3,10).Value = 297.66666
$ws.Cells.Item(3,11).Value = 3075
$ws.Cells.Item(3,12).Value = 297.66666
$ws.Cells.Item(3,13).Value = -2959
$ws.Cells.Item(3,14).Value = -529.66666

# row 4
$ws.Cells.Item(4,8).Value = 4750
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(4,13).Value = ""

# row 96
$ws.Cells.Item(96,8).Value = 9000
$ws.Cells.Item(96,10).Value = 9000
$ws.Cells.Item(96,12).Value = 9000
$ws.Cells.Item(96,14).Value = -14492

# row 102
$ws.Cells.Item(102,8).Value = 32279.922
$ws.Cells.Item(102,9).Value = 50933.24
$ws.Cells.Item(102,10).Value = 9237.588
$ws.Cells.Item(102,11).Value = 50933.24
$ws.Cells.Item(102,12).Value = 9237.588
$ws.Cells.Item(102,13).Value = -49311.24
$ws.Cells.Item(102,14).Value = -12481.588

# row 107
$ws.Cells.Item(107,8).Value = 469
$ws.Cells.Item(107,9).Value = 205.33333
$ws.Cells.Item(107,10).Value = 600.8333
$ws.Cells.Item(107,11).Value = 205.33333
$ws.Cells.Item(107,12).Value = 600.8333
$ws.Cells.Item(107,13).Value = 1714.66667
$ws.Cells.Item(107,14).Value = -4440.8333

# row 126
$ws.Cells.Item(126,8).Value = 56151
$ws.Cells.Item(126,9).Value = 59467.312
$ws.Cells.Item(126,10).Value = 3090
$ws.Cells.Item(126,11).Value = 178401.936
$ws.Cells.Item(126,12).Value = 9270
$ws.Cells.Item(126,13).Value = -175931.936
$ws.Cells.Item(126,14).Value = -14210

# row 132
$ws.Cells.Item(132,8).Value = 3339.4546
$ws.Cells.Item(132,9).Value = 3308
$ws.Cells.Item(132,10).Value = 4000
$ws.Cells.Item(132,11).Value = 9924
$ws.Cells.Item(132,12).Value = 12000
$ws.Cells.Item(132,13).Value = -7394
$ws.Cells.Item(132,14).Value = -17060

# --- LTW ---
$ws = $wb.Worksheets.Item(7)
# row 46
$ws.Cells.Item(46,8).Value = 3358.1428
$ws.Cells.Item(46,9).Value = 1084.75
$ws.Cells.Item(46,10).Value = 3651.484
$ws.Cells.Item(46,11).Value = 1084.75
$ws.Cells.Item(46,12).Value = 3651.484
$ws.Cells.Item(46,13).Value = -896.75
$ws.Cells.Item(46,14).Value = -4027.484

# row 69
$ws.Cells.Item(69,8).Value = 61259.2
$ws.Cells.Item(69,10).Value = 61824.25
$ws.Cells.Item(69,12).Value = 61824.25
$ws.Cells.Item(69,14).Value = -63446.25

# row 72
$ws.Cells.Item(72,8).Value = 61259.2
$ws.Cells.Item(72,10).Value = 61824.25
$ws.Cells.Item(72,12).Value = 185472.75
$ws.Cells.Item(72,14).Value = -193584.75

# row 100
$ws.Cells.Item(100,8).Value = 4122.3335
$ws.Cells.Item(100,9).Value = 3314.7144
$ws.Cells.Item(100,10).Value = 5253
$ws.Cells.Item(100,11).Value = 3314.7144
$ws.Cells.Item(100,12).Value = 5253
$ws.Cells.Item(100,13).Value = -2773.7144
$ws.Cells.Item(100,14).Value = -6335

# row 132
$ws.Cells.Item(132,8).Value = 2515
$ws.Cells.Item(132,9).Value = 2091.303
$ws.Cells.Item(132,10).Value = 4512.4287
$ws.Cells.Item(132,11).Value = 6273.909
$ws.Cells.Item(132,12).Value = 13537.2861
$ws.Cells.Item(132,13).Value = -3743.909
$ws.Cells.Item(132,14).Value = -18597.2861

# --- WVR ---
$ws = $wb.Worksheets.Item(8)
# row 122
$ws.Cells.Item(122,8).Value = 6425.857
$ws.Cells.Item(122,9).Value = 6197.4
$ws.Cells.Item(122,10).Value = 6997
$ws.Cells.Item(122,11).Value = 18592.2
$ws.Cells.Item(122,12).Value = 20991
$ws.Cells.Item(122,13).Value = -16142.2
$ws.Cells.Item(122,14).Value = -25891

# row 125
$ws.Cells.Item(125,8).Value = 63996.2
$ws.Cells.Item(125,10).Value = 63996.2
$ws.Cells.Item(125,12).Value = 63996.2
$ws.Cells.Item(125,14).Value = -73836.2

# row 139
$ws.Cells.Item(139,8).Value = 78326.43
$ws.Cells.Item(139,10).Value = 85547.5
$ws.Cells.Item(139,12).Value = 85547.5
$ws.Cells.Item(139,14).Value = -95827.5

# row 141
$ws.Cells.Item(141,8).Value = 161499.6
$ws.Cells.Item(141,10).Value = 161499.6
$ws.Cells.Item(141,12).Value = 161499.6
